$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 941
$ws.Range("F3").Value = 1033
$ws.Range("F4").Value = 826
$ws.Range("F5").Value = 894
$ws.Range("F7").Value = 728
$ws.Range("F8").Value = 172
$ws.Range("F9").Value = 1333
$ws.Range("F10").Value = 765
$ws.Range("F11").Value = 432
$ws.Range("F12").Value = 575
$ws.Range("F14").Value = 74
$ws.Range("F15").Value = 74
$ws.Range("F16").Value = 1317
$ws.Range("F17").Value = 153
$ws.Range("F18").Value = 21
$ws.Range("F19").Value = 438
$ws.Range("F20").Value = 10
$ws.Range("F23").Value = 610
$ws.Range("F24").Value = 170
$ws.Range("F25").Value = 673
$ws.Range("F26").Value = 42
$ws.Range("F27").Value = 1159
$ws.Range("F28").Value = 17
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 121
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 941
$ws.Range("F5").Value = 1033
$ws.Range("F6").Value = 826
$ws.Range("F7").Value = 894
$ws.Range("F9").Value = 728
$ws.Range("F10").Value = 172
$ws.Range("F11").Value = 1333
$ws.Range("F12").Value = 765
$ws.Range("F15").Value = 432
$ws.Range("F16").Value = 575
$ws.Range("F19").Value = 74
$ws.Range("F20").Value = 74
$ws.Range("F21").Value = 1317
$ws.Range("F23").Value = 153
$ws.Range("F24").Value = 21
$ws.Range("F25").Value = 438
$ws.Range("F26").Value = 10
$ws.Range("F31").Value = 610
$ws.Range("F34").Value = 121
$ws.Range("F35").Value = 121
$ws.Range("F36").Value = 170
$ws.Range("F37").Value = 673
$ws.Range("F38").Value = 42
$ws.Range("F39").Value = 1159
$ws.Range("F40").Value = 18
